$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @([double]"0.06282870471477509",[double]"0.9828430414199829",[double]"0.01085071451961994",[double]"0.9979230165481567"),
    @([double]"0.01039525680243969",[double]"0.9982308149337769",[double]"0.00875630509108305",[double]"0.9979230165481567"),
    @([double]"0.007095101289451122",[double]"0.998457133769989",[double]"0.005982016678899527",[double]"0.99806147813797"),
    @([double]"0.003250523470342159",[double]"0.9989919662475586",[double]"0.001100260997191072",[double]"0.9993076920509338"),
    @([double]"0.002281023422256112",[double]"0.9995062947273254",[double]"0.0003513776755426079",[double]"0.9998615384101868"),
    @([double]"0.001272676978260279",[double]"0.9996502995491028",[double]"0.0008334434824064374",[double]"0.9997230768203735"),
    @([double]"0.001134548685513437",[double]"0.9997119903564453",[double]"0.0002371456648688763",[double]"1"),
    @([double]"0.0007800249150022864",[double]"0.9997736811637878",[double]"0.0003103387134615332",[double]"0.9998615384101868"),
    @([double]"0.00178591173607856",[double]"0.9996502995491028",[double]"0.000146520949783735",[double]"1"),
    @([double]"0.0008525385637767613",[double]"0.9996914267539978",[double]"0.0001066033728420734",[double]"1"),
    @([double]"0.0007122239330783486",[double]"0.9998354315757751",[double]"0.0002789738646242768",[double]"0.9998615384101868"),
    @([double]"0.0005686317454092205",[double]"0.9998354315757751",[double]"3.645956212494639E-06",[double]"1"),
    @([double]"0.001130313496105373",[double]"0.9997325539588928",[double]"2.455276990076527E-05",[double]"1"),
    @([double]"0.0003081322356592864",[double]"0.9998765587806702",[double]"5.720556146115996E-06",[double]"1"),
    @([double]"0.0002759694471023977",[double]"0.9998765587806702",[double]"0.000296323501970619",[double]"0.9997230768203735"),
    @([double]"0.0004831361875403672",[double]"0.9999383091926575",[double]"2.068367393803783E-06",[double]"1"),
    @([double]"0.0005349786952137947",[double]"0.9998354315757751",[double]"0.000638014986179769",[double]"0.9997230768203735"),
    @([double]"0.0005231896648183465",[double]"0.9998148679733276",[double]"0.00359085900709033",[double]"0.9997230768203735"),
    @([double]"0.0008616555132903159",[double]"0.9998354315757751",[double]"0.002779387170448899",[double]"0.9997230768203735"),
    @([double]"0.0006011150544509292",[double]"0.9998971223831177",[double]"1.653091931075323E-05",[double]"1"),
    @([double]"0.0002357810008106753",[double]"0.999958872795105",[double]"2.828609467542265E-06",[double]"1"),
    @([double]"2.041883453784976E-05",[double]"1",[double]"1.050290211423999E-05",[double]"1"),
    @([double]"0.001388830598443747",[double]"0.9998148679733276",[double]"2.844104983523721E-06",[double]"1"),
    @([double]"0.000322646606946364",[double]"0.9998765587806702",[double]"2.761302084763884E-06",[double]"1"),
    @([double]"0.000334899581503123",[double]"0.9999383091926575",[double]"7.479215855710208E-05",[double]"1"),
    @([double]"0.0001319392613368109",[double]"0.999958872795105",[double]"4.18257259298116E-05",[double]"1"),
    @([double]"6.819709960836917E-05",[double]"0.999958872795105",[double]"0.0001725208421703428",[double]"0.9998615384101868"),
    @([double]"0.00052812130888924",[double]"0.9999176859855652",[double]"2.240926960439538E-06",[double]"1"),
    @([double]"1.930411781358998E-05",[double]"1",[double]"7.000030018389225E-05",[double]"1"),
    @([double]"2.858807238226291E-05",[double]"0.9999794363975525",[double]"6.076084559936135E-07",[double]"1"),
    @([double]"0.0001593780471011996",[double]"0.999958872795105",[double]"0.0004123479302506894",[double]"0.9998615384101868"),
    @([double]"0.0009153003338724375",[double]"0.9998971223831177",[double]"0.0003424502210691571",[double]"0.9998615384101868"),
    @([double]"0.0004767657956108451",[double]"0.9999176859855652",[double]"1.358641839033226E-06",[double]"1"),
    @([double]"0.0001063353411154822",[double]"0.999958872795105",[double]"1.23450035971473E-07",[double]"1"),
    @([double]"1.737596721795853E-05",[double]"1",[double]"1.276242187486787E-06",[double]"1"),
    @([double]"0.0003067961079068482",[double]"0.9999383091926575",[double]"1.210474920299021E-06",[double]"1"),
    @([double]"0.0004406876105349511",[double]"0.9998971223831177",[double]"9.759023811284351E-08",[double]"1"),
    @([double]"0.0003309159365016967",[double]"0.999958872795105",[double]"4.71189650852466E-07",[double]"1"),
    @([double]"0.0005877171061001718",[double]"0.9999176859855652",[double]"2.990876879493953E-08",[double]"1"),
    @([double]"8.043642083066516E-06",[double]"1",[double]"9.375569653968796E-09",[double]"1"),
    @([double]"0.0002470318868290633",[double]"0.9999176859855652",[double]"8.941028681874741E-06",[double]"1"),
    @([double]"0.0002659836027305573",[double]"0.9998971223831177",[double]"0.0009957951260730624",[double]"0.9997230768203735"),
    @([double]"0.0001214111834997311",[double]"0.9999794363975525",[double]"1.445577709091594E-05",[double]"1"),
    @([double]"0.0003123309288639575",[double]"0.999958872795105",[double]"1.270274879061617E-05",[double]"1"),
    @([double]"0.0003110425022896379",[double]"0.9999383091926575",[double]"1.616396048120805E-06",[double]"1"),
    @([double]"0.0004860930785071105",[double]"0.9998765587806702",[double]"1.115945224228199E-06",[double]"1"),
    @([double]"2.635229975567199E-06",[double]"1",[double]"9.829424243434914E-07",[double]"1"),
    @([double]"0.000124282407341525",[double]"0.999958872795105",[double]"1.333587874796649E-06",[double]"1"),
    @([double]"6.188482802826911E-05",[double]"0.9999794363975525",[double]"1.006882222043259E-08",[double]"1"),
    @([double]"0.0001330009836237878",[double]"0.999958872795105",[double]"8.201206583180465E-06",[double]"1")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row++
}
